# Applies the "Add files via upload" edit:
#  - Row for "2506 - KATIA FERREIRA DE BARROS" becomes
#    "2515 - MARCO FREIRE (ÁREA EXTERNA)"
#  - A new row "2516 - JOSÉ CARLOS MORAES ABREU FILHO" is appended to the
#    Graziele/SP block (the old "2314 - LILIAN STUHLBERGER" row is removed,
#    so the net row count is unchanged; the rows in between shift up by one)
#  - "Alana" is renamed to "Roberto"
#  - Row for "2405 - OTÁVIO UMADA" becomes "2514 - FELIPE HESS BORGES"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. 2506 - KATIA FERREIRA DE BARROS -> 2515 - MARCO FREIRE (ÁREA EXTERNA)  (row 18)
$ws.Cells.Item(18, 1).Value = 2515
$ws.Cells.Item(18, 2).Value = "2515 - MARCO FREIRE (ÁREA EXTERNA)"

# 2. The "2314 - LILIAN STUHLBERGER" row (row 6) is removed and a new
#    "2516 - JOSÉ CARLOS MORAES ABREU FILHO" entry is added at the end of
#    the Graziele/SP block (row 12) -- net effect: rows 7..12 shift up into
#    6..11, and row 12 becomes the new entry.
$ws.Cells.Item(6, 1).Value = 2407
$ws.Cells.Item(6, 2).Value = "2407 - SUN MORITZ ADMINISTRADORA"
$ws.Cells.Item(7, 1).Value = 2412
$ws.Cells.Item(7, 2).Value = "2412 - GABRIEL ACURCIO V. S. DE CARVALHO"
$ws.Cells.Item(8, 1).Value = 2413
$ws.Cells.Item(8, 2).Value = "2413 - FERNANDO VASCONCELLOS"
$ws.Cells.Item(9, 1).Value = 2503
$ws.Cells.Item(9, 2).Value = "2503 - KAENA PARTICIPAÇÕES LTDA"
$ws.Cells.Item(10, 1).Value = 2511
$ws.Cells.Item(10, 2).Value = "2511 - 1807 PARTICIPAÇÕES LTDA"
$ws.Cells.Item(11, 1).Value = 2512
$ws.Cells.Item(11, 2).Value = "2512 - ROBERTO KLABIN MARTINS XAVIER"
$ws.Cells.Item(12, 1).Value = 2516
$ws.Cells.Item(12, 2).Value = "2516 - JOSÉ CARLOS MORAES ABREU FILHO"

# 3. Alana -> Roberto (row 13)
$ws.Cells.Item(13, 3).Value = "Roberto"

# 4. 2405 - OTÁVIO UMADA -> 2514 - FELIPE HESS BORGES (row 3)
$ws.Cells.Item(3, 1).Value = 2514
$ws.Cells.Item(3, 2).Value = "2514 - FELIPE HESS BORGES"

# Update the selection to match the saved view state.
$ws.Range("F5").Select()
